# Update the "展览" (Exhibition) and "全部类型" (All types) sheets to reflect
# the refreshed scrape of the convention listing:
#  - three events that have passed / were superseded are removed
#    (rows that originally held "新余·新次元动漫游戏嘉年华",
#     "江西·ShiningStaR数字互娱嘉年华", and its companion meet-and-greet row)
#  - the remaining events' "想去人数" (interest count) values are refreshed

$wb = $excel.ActiveWorkbook

# New "想去人数" (column F) values for the 17 surviving events, in final row order
# (rows 2-18 after the three obsolete rows have been removed).
$newCounts = @{
    2  = 1541
    3  = 122
    4  = 397
    5  = 280
    6  = 35
    7  = 137
    8  = 48
    9  = 468
    10 = 1298
    11 = 356
    12 = 88
    13 = 153
    14 = 107
    15 = 148
    16 = 87
    17 = 135
    18 = 122
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Remove the three obsolete rows. Delete from the bottom up so earlier
    # deletions don't shift the row numbers of rows still to be removed.
    $ws.Rows.Item(8).Delete()
    $ws.Rows.Item(3).Delete()
    $ws.Rows.Item(2).Delete()

    # Refresh the "想去人数" counts for the remaining data rows.
    foreach ($r in $newCounts.Keys) {
        $ws.Cells.Item($r, 6).Value = $newCounts[$r]
    }
}
